$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B/C values between row 2 and row 3
$ws.Range("B2").Value = "12f3026e-afac-400e-a3c9-7ae04c185ad7"
$ws.Range("C2").Value = 1

$ws.Range("B3").Value = "01849493-c3a8-49c2-ab66-2f6c6606e6d1"
$ws.Range("C3").Value = 3

# Give row 3 a custom (slightly taller) height
$ws.Rows.Item(3).RowHeight = 15.6

# Update the selection on the sheet to A1:C2 (no explicit active cell, defaults to top-left)
$ws.Range("A1:C2").Select()
